# Update the Expense sheet with new loan records.
# Rows 2-4 are overwritten with new values, and a new row 5 is appended.
# Date/time values must stay plain text (shared strings), not be coerced
# into Excel date serials, so they are written as explicit strings.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = 11
$ws.Cells.Item(2, 2).Value = "2022-01-19 20:07:54.443948"
$ws.Cells.Item(2, 3).Value = 10000
$ws.Cells.Item(2, 4).Value = "USD"
$ws.Cells.Item(2, 5).Value = "Jack"
$ws.Cells.Item(2, 6).Value = "John"

# Row 3
$ws.Cells.Item(3, 1).Value = 13
$ws.Cells.Item(3, 2).Value = "2022-01-19 21:51:18.819872"
$ws.Cells.Item(3, 3).Value = 58895544
$ws.Cells.Item(3, 4).Value = "UAH"
$ws.Cells.Item(3, 5).Value = "Jack"
$ws.Cells.Item(3, 6).Value = "John"

# Row 4
$ws.Cells.Item(4, 1).Value = 14
$ws.Cells.Item(4, 2).Value = "2022-01-19 21:51:44.560407"
$ws.Cells.Item(4, 3).Value = 10000
$ws.Cells.Item(4, 4).Value = "USD"
$ws.Cells.Item(4, 5).Value = "Jack"
$ws.Cells.Item(4, 6).Value = "John"

# Row 5 (new)
$ws.Cells.Item(5, 1).Value = 15
$ws.Cells.Item(5, 2).Value = "2022-01-19 21:52:56.170675"
$ws.Cells.Item(5, 3).Value = 2375821
$ws.Cells.Item(5, 4).Value = "USD"
$ws.Cells.Item(5, 5).Value = "John"
$ws.Cells.Item(5, 6).Value = "Jack"
